$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cells carrying the correct style + shared-string value for PRESENT/ABSENT
$presentRef = $ws.Range("G7")
$absentRef = $ws.Range("N7")

# Step 1: fill the whole Q7:S79 block with PRESENT (the default for this update)
$fillRange = $ws.Range("Q7:S79")
$fillRange.Value = $presentRef.Value2
$presentRef.Copy()
$fillRange.PasteSpecial(-4122)

# Step 2: correct the handful of cells that are actually ABSENT
$absentCells = "S28,Q32,S33,Q37,R37,Q38,R50,Q61,Q63,Q67,R77"
foreach ($addr in $absentCells.Split(",")) {
    $cell = $ws.Range($addr)
    $cell.Value = $absentRef.Value2
    $absentRef.Copy()
    $cell.PasteSpecial(-4122)
}

# Step 3: update the Total Present (F) and, where needed, Total Absence (E) summary cells
$ws.Range("F7").Value = 11
$ws.Range("F8").Value = 12
$ws.Range("F9").Value = 13
$ws.Range("F10").Value = 13
$ws.Range("F11").Value = 11
$ws.Range("F12").Value = 11
$ws.Range("F13").Value = 10
$ws.Range("F14").Value = 11
$ws.Range("F15").Value = 10
$ws.Range("F16").Value = 13
$ws.Range("F17").Value = 9
$ws.Range("F18").Value = 11
$ws.Range("F19").Value = 9
$ws.Range("F20").Value = 11
$ws.Range("F21").Value = 8
$ws.Range("F22").Value = 11
$ws.Range("F23").Value = 9
$ws.Range("F24").Value = 10
$ws.Range("F25").Value = 9
$ws.Range("F26").Value = 11
$ws.Range("F27").Value = 11
$ws.Range("F28").Value = 9
$ws.Range("E28").Value = 4
$ws.Range("F29").Value = 13
$ws.Range("F30").Value = 12
$ws.Range("F31").Value = 12
$ws.Range("F32").Value = 11
$ws.Range("E32").Value = 2
$ws.Range("F33").Value = 10
$ws.Range("E33").Value = 3
$ws.Range("F34").Value = 11
$ws.Range("F35").Value = 10
$ws.Range("F36").Value = 12
$ws.Range("F37").Value = 8
$ws.Range("E37").Value = 5
$ws.Range("F38").Value = 9
$ws.Range("E38").Value = 4
$ws.Range("F39").Value = 10
$ws.Range("F40").Value = 12
$ws.Range("F41").Value = 11
$ws.Range("F42").Value = 11
$ws.Range("F43").Value = 12
$ws.Range("F44").Value = 12
$ws.Range("F45").Value = 12
$ws.Range("F46").Value = 12
$ws.Range("F47").Value = 11
$ws.Range("F48").Value = 10
$ws.Range("F49").Value = 11
$ws.Range("F50").Value = 11
$ws.Range("E50").Value = 2
$ws.Range("F51").Value = 12
$ws.Range("F52").Value = 10
$ws.Range("F53").Value = 11
$ws.Range("F54").Value = 9
$ws.Range("F55").Value = 12
$ws.Range("F56").Value = 10
$ws.Range("F57").Value = 11
$ws.Range("F58").Value = 13
$ws.Range("F59").Value = 11
$ws.Range("F60").Value = 11
$ws.Range("F61").Value = 10
$ws.Range("E61").Value = 3
$ws.Range("F62").Value = 12
$ws.Range("F63").Value = 11
$ws.Range("E63").Value = 2
$ws.Range("F64").Value = 11
$ws.Range("F65").Value = 10
$ws.Range("F66").Value = 11
$ws.Range("F67").Value = 11
$ws.Range("E67").Value = 2
$ws.Range("F68").Value = 12
$ws.Range("F69").Value = 12
$ws.Range("F70").Value = 13
$ws.Range("F71").Value = 10
$ws.Range("F72").Value = 11
$ws.Range("F73").Value = 11
$ws.Range("F74").Value = 11
$ws.Range("F75").Value = 12
$ws.Range("F76").Value = 10
$ws.Range("F77").Value = 11
$ws.Range("E77").Value = 2
$ws.Range("F78").Value = 11
$ws.Range("F79").Value = 9

# Row 80: the PRESENT-style template bleeds one cell down into the still-empty row (Q80)
$presentRef.Copy()
$ws.Range("Q80").PasteSpecial(-4122)
$ws.Range("Q80").ClearContents()

$excel.CutCopyMode = $false